$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# Slide 1: remove the stale workshop date/time from the subtitle
# ---------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$subtitle = $s1.Shapes.Item(2)
$subtitle.TextFrame.TextRange.Text = ""

# ---------------------------------------------------------------
# Slide 13 ("Loops"): fix the sample code block
#   - "$ file_list=`ls`"   -> the `ls` run keeps its trailing back-tick
#   - continuation prompts "$ echo $file" / "$ done" become "> echo
#     $file" / "> done" (bash PS2-style prompt)
# ---------------------------------------------------------------
$s13 = $p.Slides.Item(13)
$code = $s13.Shapes.Item(2)
$trCode = $code.TextFrame.TextRange

# Edit from the end of the range backwards so earlier offsets stay valid.
$trCode.Characters(146, 2).Text = "> "   # "$ " (of "$ done") -> "> "
$trCode.Characters(133, 2).Text = "> "   # "$ " (of "$ echo $file") -> "> "
$trCode.Characters(102, 1).Text = ""     # trailing "`" run folds into "ls`"
$trCode.Characters(100, 2).Text = "ls``" # "ls" -> "ls`"
